$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the two time-range values in column C
$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"

# Update the selection on the sheet: active cell stays C11, but selected range extends to C10:C11
$ws.Range("C10:C11").Select()
$ws.Range("C11").Activate()
